# 学生信息表 - append 7 more student rows (0010002..0010008), mirroring the
# existing header/row-2 formatting: column A ids are stored as text (quote
# prefix, like the existing "0010001"), column D is a date using the same
# m/d/yy date format as the existing date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 3; Id = "0010002"; Name = "张四";   Gender = "女"; Date = 40911 },
    @{ Row = 4; Id = "0010003"; Name = "流量";   Gender = "女"; Date = 40912 },
    @{ Row = 5; Id = "0010004"; Name = "驴";     Gender = "女"; Date = 40913 },
    @{ Row = 6; Id = "0010005"; Name = "驴2";    Gender = "女"; Date = 40914 },
    @{ Row = 7; Id = "0010006"; Name = "李老师"; Gender = "女"; Date = 40915 },
    @{ Row = 8; Id = "0010007"; Name = "达到";   Gender = "女"; Date = 40916 },
    @{ Row = 9; Id = "0010008"; Name = "张三2";  Gender = "女"; Date = 40917 }
)

foreach ($r in $rows) {
    $n = $r.Row
    # Leading apostrophe forces text storage (same quote-prefix style used by
    # the existing "0010001" id in A2) instead of Excel auto-coercing it to
    # a number.
    $ws.Range("A$n").Value = "'" + $r.Id
    $ws.Range("B$n").Value = $r.Name
    $ws.Range("C$n").Value = $r.Gender
    $ws.Range("D$n").Value = $r.Date
    $ws.Range("D$n").NumberFormat = "m/d/yy"
}

[void]$ws.Range("D8:D9").Select()
